# Generate Report for Handoff
#
# The "b.md" row (row 3) in each of the three report sheets moves from
# "Handed back: in sync with en-US" to "Ready for handoff": new handoff
# file names/timestamps are recorded and a version-mismatch error message
# is attached. The Overview sheet mirrors the new status + timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: b.md row (row 3) ------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-27 12:35:17"

# --- zh-cn sheet: b.md row (row 3) ----------------------------------------
# (the leading "'" forces the "False" literal to stay plain text instead of
# being auto-coerced into an actual Boolean by COM's type inference; resetting
# the style back to "Normal" afterwards drops the quote-prefix marker again
# since this cell carries no special formatting)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-27 12:35:12"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a78039e5a46df065785e576854a76c5a0aac7894/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18e132661963b1a5623da574b6bb649747cb1800/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666666

# --- de-de sheet: b.md row (row 3) ----------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-27 12:35:17"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a78039e5a46df065785e576854a76c5a0aac7894/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18e132661963b1a5623da574b6bb649747cb1800/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.16666666666666
